# Trade #98 closed at 2026-02-18 00:33:58 - unknown UNKNOWN +0.000%
#
# This script applies the edits described by the commit diff:
#   - Summary sheet: bump Total Trades / Win Rate %
#   - Strategy Status sheet: bump MarketMaking Trades / Win Rate %
#   - All Trades sheet: close trade #127 (row 128) and append 2 new open
#     trades (rows 157/158)
#   - momentum sheet: append the new momentum trade (row 40)
#   - HighProbConvergence sheet: append the new HighProbConvergence trade (row 21)
#   - MarketMaking sheet: close trade #127 (row 48)

$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-TextCell($ws, $addr, $val) {
    # Force a literal text value, even when it looks like a date/number,
    # by pre-formatting the cell as Text before writing it.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
Set-Cell $wsSummary "B6" 126
Set-Cell $wsSummary "B9" 46.83

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
Set-Cell $wsStatus "D6" 46
Set-Cell $wsStatus "G6" 45.65

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close existing trade #127 (row 128)
Set-Cell $wsAll "G128" 0.98
Set-Cell $wsAll "H128" "CLOSED"
Set-Cell $wsAll "K128" 99.54000000000001
Set-Cell $wsAll "L128" "early_exit"
Set-Cell $wsAll "M128" 0.13

# Append new trade #156 -> row 157 (momentum, still OPEN)
Set-Cell $wsAll "A157" 156
Set-TextCell $wsAll "B157" "2026-02-18"
Set-Cell $wsAll "C157" "00:33:52"
Set-Cell $wsAll "D157" "momentum"
Set-Cell $wsAll "E157" "UP"
Set-Cell $wsAll "F157" 0.98
Set-Cell $wsAll "H157" "OPEN"
Set-Cell $wsAll "I157" 0
Set-Cell $wsAll "J157" 0
Set-Cell $wsAll "K157" 99.23374292899115
Set-Cell $wsAll "M157" 0
Set-Cell $wsAll "N157" 0
Set-Cell $wsAll "O157" 0
Set-Cell $wsAll "P157" 0.9
Set-Cell $wsAll "Q157" "Upward momentum: 1.980% over 10 samples"

# Append new trade #157 -> row 158 (HighProbConvergence, still OPEN)
Set-Cell $wsAll "A158" 157
Set-TextCell $wsAll "B158" "2026-02-18"
Set-Cell $wsAll "C158" "00:33:52"
Set-Cell $wsAll "D158" "HighProbConvergence"
Set-Cell $wsAll "E158" "DOWN"
Set-Cell $wsAll "F158" 0.01
Set-Cell $wsAll "H158" "OPEN"
Set-Cell $wsAll "I158" 0
Set-Cell $wsAll "J158" 0
Set-Cell $wsAll "K158" 100.4130057263667
Set-Cell $wsAll "M158" 0
Set-Cell $wsAll "N158" 0
Set-Cell $wsAll "O158" 0
Set-Cell $wsAll "P158" 0.95
Set-Cell $wsAll "Q158" "Mean reversion DOWN: price 1.68% above mean (z=2.38)"

# ---------------------------------------------------------------------
# momentum sheet - append new trade #156 -> row 40 (OPEN)
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
Set-Cell $wsMomentum "A40" 156
Set-TextCell $wsMomentum "B40" "2026-02-18"
Set-Cell $wsMomentum "C40" "00:33:52"
Set-Cell $wsMomentum "D40" "momentum"
Set-Cell $wsMomentum "E40" "UP"
Set-Cell $wsMomentum "F40" 0.98
Set-Cell $wsMomentum "H40" "OPEN"
Set-Cell $wsMomentum "I40" 0
Set-Cell $wsMomentum "J40" 0
Set-Cell $wsMomentum "K40" 99.23374292899115
Set-Cell $wsMomentum "L40" 0
Set-Cell $wsMomentum "M40" 0
Set-Cell $wsMomentum "N40" 0.9
Set-Cell $wsMomentum "O40" "Upward momentum: 1.980% over 10 samples"
Set-Cell $wsMomentum "Q40" 0

# ---------------------------------------------------------------------
# HighProbConvergence sheet - append new trade #157 -> row 21 (OPEN)
# ---------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
Set-Cell $wsHPC "A21" 157
Set-TextCell $wsHPC "B21" "2026-02-18"
Set-Cell $wsHPC "C21" "00:33:52"
Set-Cell $wsHPC "D21" "HighProbConvergence"
Set-Cell $wsHPC "E21" "DOWN"
Set-Cell $wsHPC "F21" 0.01
Set-Cell $wsHPC "H21" "OPEN"
Set-Cell $wsHPC "I21" 0
Set-Cell $wsHPC "J21" 0
Set-Cell $wsHPC "K21" 100.4130057263667
Set-Cell $wsHPC "L21" 0
Set-Cell $wsHPC "M21" 0
Set-Cell $wsHPC "N21" 0.95
Set-Cell $wsHPC "O21" "Mean reversion DOWN: price 1.68% above mean (z=2.38)"
Set-Cell $wsHPC "Q21" 0

# ---------------------------------------------------------------------
# MarketMaking sheet - close existing trade #127 (row 48)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
Set-Cell $wsMM "G48" 0.98
Set-Cell $wsMM "H48" "CLOSED"
Set-Cell $wsMM "K48" 99.54000000000001
Set-Cell $wsMM "P48" "early_exit"
Set-Cell $wsMM "Q48" 0.13

Write-Output "Edits applied."
